$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-08-29 09:10:25"
$wsZhCn.Range("H2").Value = "2016-08-29 09:10:20"
$wsZhCn.Range("K2").Value = "2016-08-29 09:10:38"
$wsDeDe.Range("H2").Value = "2016-08-29 09:10:25"
$wsDeDe.Range("K2").Value = "2016-08-29 09:10:45"
